# Append the new Adafruit IO reading as row 55 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# Force column C to text so a numeric-looking reading ("25") is stored
# as a string, matching the other rows in this feed-log sheet.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
